$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Supplier" column (K): a header cell plus a supplier id of 2
# for every existing product row, so the importer test fixture has a
# supplier value to feed into the import.
$ws.Range("K1").Value = "Supplier"
$ws.Range("K2:K7").Value = 2

# The edited workbook keeps row 2 slightly shorter than the rest.
$ws.Rows.Item(2).RowHeight = 15

# Leave the selection on the newly added cell, matching the saved view
# state of the edited workbook.
$ws.Range("K7").Select()
